# Update course details for ACE AVIATION AEROSPACE ACADEMY workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# department: "ACE AVIATION" -> "AVIATION"
$ws.Range("C2").Value = "AVIATION"

# promotionValidity: clear the "Promotion valid until  31th Dec 2021" text
$ws.Range("R2").Value = ""
